$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing data
# (old A..E) one column to the right (new B..F).
$ws.Columns("A").Insert()

# Populate the new column A with the method names (header + 8 methods).
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Update the header row (previously Var1_1..Var1_5, shifted to B1..F1).
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Recompute the best-fit column widths: column A now holds the longer
# method names, B/C hold short numeric codes, D/F keep their previous
# (already best-fit) width inherited from the shift caused by the insert.
$ws.Columns("A").ColumnWidth = 12.5
$ws.Columns("B:C").ColumnWidth = 2.33
